$d = $word.ActiveDocument

# Locate the date text "Күні: 13.10.2023ж" and edit it to "Күні: 16.02.2024ж"
# by retyping the day, month, and last digit of the year in place, mirroring
# how a human editor would select and retype each piece (which is what
# produced the run-split pattern in the target document).

# --- Step 1: day "3" -> "6" (13 -> 16) ---
$find1 = $d.Content
$found1 = $find1.Find.Execute("Күні: 1", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) { throw "Could not find 'Kuni: 1' anchor" }
$dayRng = $d.Range($find1.End, $find1.End + 1)
$dayRng.Text = "6"

# --- Step 2: month "10" -> "02" ---
$find2 = $d.Content
$found2 = $find2.Find.Execute("16.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) { throw "Could not find '16.' anchor" }
$monthRng = $d.Range($find2.End, $find2.End + 2)
$monthRng.Text = "02"

# --- Step 3: year last digit "3" -> "4" (2023 -> 2024) ---
$find3 = $d.Content
$found3 = $find3.Find.Execute("16.02.202", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found3) { throw "Could not find '16.02.202' anchor" }
$yearRng = $d.Range($find3.End, $find3.End + 1)
$yearRng.Text = "4"
